$wb = $excel.ActiveWorkbook

# Rows updated in each localization report sheet
$rows = @(7, 8, 11, 12, 13, 14)

# --- Overview sheet: update "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-13 18:25:11"
}

# --- zh-cn sheet: update "Priority" (E) and "Latest Handoff Datetime" (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-13 18:24:59"
}

# --- de-de sheet: update "Priority" (E) and "Latest Handoff Datetime" (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-13 18:25:11"
}
